# Update "想去人数" (interested-people count) values in column F
# for the 展览 (sheet1), 演出 (sheet2) and 全部类型 (sheet4) worksheets,
# reflecting a refreshed data pull (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 5115
$ws1.Range("F4").Value  = 8
$ws1.Range("F5").Value  = 7420
$ws1.Range("F8").Value  = 100
$ws1.Range("F10").Value = 67
$ws1.Range("F11").Value = 25
$ws1.Range("F12").Value = 4309
$ws1.Range("F13").Value = 1750
$ws1.Range("F14").Value = 100
$ws1.Range("F15").Value = 104
$ws1.Range("F16").Value = 2910
$ws1.Range("F17").Value = 579
$ws1.Range("F18").Value = 566
$ws1.Range("F20").Value = 495
$ws1.Range("F21").Value = 432
$ws1.Range("F22").Value = 453
$ws1.Range("F23").Value = 305
$ws1.Range("F24").Value = 98
$ws1.Range("F25").Value = 1688
$ws1.Range("F26").Value = 1182
$ws1.Range("F27").Value = 90
$ws1.Range("F28").Value = 1376
$ws1.Range("F32").Value = 514
$ws1.Range("F34").Value = 59
$ws1.Range("F36").Value = 61
$ws1.Range("F37").Value = 2865
$ws1.Range("F38").Value = 700
$ws1.Range("F39").Value = 16
$ws1.Range("F40").Value = 59
$ws1.Range("F42").Value = 16

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 7
$ws2.Range("F3").Value = 9

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 5115
$ws4.Range("F4").Value  = 8
$ws4.Range("F5").Value  = 7420
$ws4.Range("F8").Value  = 100
$ws4.Range("F10").Value = 67
$ws4.Range("F11").Value = 25
$ws4.Range("F12").Value = 4308
$ws4.Range("F13").Value = 1750
$ws4.Range("F14").Value = 100
$ws4.Range("F15").Value = 104
$ws4.Range("F16").Value = 2910
$ws4.Range("F17").Value = 579
$ws4.Range("F18").Value = 566
$ws4.Range("F20").Value = 496
$ws4.Range("F21").Value = 432
$ws4.Range("F22").Value = 453
$ws4.Range("F23").Value = 7
$ws4.Range("F24").Value = 305
$ws4.Range("F25").Value = 98
$ws4.Range("F26").Value = 1688
$ws4.Range("F27").Value = 1182
$ws4.Range("F28").Value = 90
$ws4.Range("F29").Value = 1376
$ws4.Range("F33").Value = 514
$ws4.Range("F35").Value = 59
$ws4.Range("F37").Value = 61
$ws4.Range("F38").Value = 2865
$ws4.Range("F39").Value = 9
$ws4.Range("F40").Value = 700
$ws4.Range("F41").Value = 16
$ws4.Range("F42").Value = 59
$ws4.Range("F44").Value = 16
